# Convention change to support multi-axle vehicles:
# add a new tire-data sheet ("Tir_430_50R38") for the Truck_430_50R38 tire,
# built from a copy of the last existing tire sheet.

$wb = $excel.ActiveWorkbook

# Duplicate the last sheet ("Tir_145_70R13") - this carries over all
# formatting, column widths, conditional formatting and data validation
# exactly as the new sheet needs them. The copy is placed right after the
# source sheet, i.e. at the end of the tab strip.
$sourceSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sourceSheet.Copy($null, $sourceSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Tir_430_50R38"

# Update the MFSwift instance name and the referenced .tir file for the
# new 430/50R38 truck tire (everything else on the sheet stays identical
# to the template it was copied from).
$newSheet.Range("H3").Value = "MFSwift_430_50R38"
$newSheet.Range("H5").Value = "which('Truck_430_50R38.tir')"

# The newly added sheet becomes the active tab, with H6 selected.
$newSheet.Activate()
[void]$newSheet.Range("H6").Select()
